$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume snapshot values.
# D-column (Price) values are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the source inlineStr cells) instead
# of auto-converting numeric-looking strings (e.g. "1.00" -> 1, "9.10" -> 9.1).

$ws.Range("D2").Value = '''45.280.89'
$ws.Range("E2").Value = '  +2.28%  '
$ws.Range("D3").Value = '''2.419.30'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''319.01'
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("D6").Value = '''102.79'
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +5.41%  '
$ws.Range("D10").Value = '''35.43'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").Value = '''0.0798'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("D13").Value = '''18.22'
$ws.Range("E13").Value = '  -3.49%  '
$ws.Range("D14").Value = '''7.01'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '''2.798.81'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = '''2.430.20'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '''0.835'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '''45.219.10'
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("D19").Value = '''12.21'
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").Value = '''6.32'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("D21").Value = '''0.0₃0920'
$ws.Range("D22").Value = '''70.49'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = '''244.44'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("E24").Value = '  -2.46%  '
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''25.63'
$ws.Range("E27").Value = '  +1.79%  '
$ws.Range("D28").Value = '''2.27'
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").Value = '''9.64'
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("D30").Value = '''49.36'
$ws.Range("E30").Value = '  +1.84%  '
$ws.Range("D31").Value = '''32.82'
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").Value = '''20.18'
$ws.Range("E32").Value = '  +7.53%  '
$ws.Range("D33").Value = '''0.125'
$ws.Range("E33").Value = '  +6.79%  '
$ws.Range("D34").Value = '''5.22'
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("D36").Value = '''0.0760'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '''1.86'
$ws.Range("E37").Value = '  -1.74%  '
$ws.Range("D38").Value = '''4.43'
$ws.Range("E38").Value = '  -0.65%  '
$ws.Range("D39").Value = '''128.39'
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '''2.27'
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").Value = '''20.47'
$ws.Range("E43").Value = '  -3.91%  '
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '''1.941.81'
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("E47").Value = '  +2.03%  '
$ws.Range("D48").Value = '''1.78'
$ws.Range("E48").Value = '  +7.32%  '
$ws.Range("D49").Value = '''9.10'
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("D50").Value = '''77.02'
$ws.Range("E50").Value = '  +4.33%  '
$ws.Range("D51").Value = '''4.79'
$ws.Range("E51").Value = '  +4.83%  '
